# Natmi following Dr Hou advice
# Rebuild the LR-pair data body (rows 2-9) on Sheet1 with the updated
# NATMI edge statistics that now also include the "sCs" sending/target
# cluster, while leaving the header row (row 1) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data body (keep header row 1 intact).
$ws.Range("A2:T5").ClearContents()

# Each entry: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric statistics columns (E..T).
$data = @(
  @("FAPs","Wnt5a","Fzd7","ECs", 3,1,10.43424333333333,31.30273,0.9711091978791583,0.9711091978791584,3,1,2.500933333333333,7.5028,0.129395502869887,0.129395502869887,26.09534696044444,234.858122644,0.1256571630011463,0.1256571630011463),
  @("FAPs","Wnt5a","Fzd7","FAPs", 3,1,10.43424333333333,31.30273,0.9711091978791583,0.9711091978791584,3,1,7.782607,23.347821,0.4026634108880828,0.4026634108880828,81.20561520570332,730.8505368513299,0.391030141962812,0.3910301419628121),
  @("FAPs","Wnt5a","Fzd7","M2", 3,1,10.43424333333333,31.30273,0.9711091978791583,0.9711091978791584,3,1,4.565521,13.696563,0.2362149673420707,0.2362149673420707,47.63775705744334,428.73981351699,0.2293905274626099,0.2293905274626099),
  @("FAPs","Wnt5a","Fzd7","sCs", 3,1,10.43424333333333,31.30273,0.9711091978791583,0.9711091978791584,3,1,4.478761333333333,13.436284,0.2317261188999595,0.2317261188999595,46.73248558392444,420.5923702553199,0.2250313654525902,0.2250313654525902),
  @("sCs","Wnt5a","Fzd7","ECs", 2,0.6666666666666666,0.310422,0.9312659999999999,0.02889080212084161,0.02889080212084161,3,1,2.500933333333333,7.5028,0.129395502869887,0.129395502869887,0.7763447271999999,6.987102544799999,0.003738339868740697,0.003738339868740697),
  @("sCs","Wnt5a","Fzd7","FAPs", 2,0.6666666666666666,0.310422,0.9312659999999999,0.02889080212084161,0.02889080212084161,3,1,7.782607,23.347821,0.4026634108880828,0.4026634108880828,2.415892430154,21.743031871386,0.01163326892527074,0.01163326892527074),
  @("sCs","Wnt5a","Fzd7","M2", 2,0.6666666666666666,0.310422,0.9312659999999999,0.02889080212084161,0.02889080212084161,3,1,4.565521,13.696563,0.2362149673420707,0.2362149673420707,1.417238159862,12.755143438758,0.006824439879460828,0.006824439879460828),
  @("sCs","Wnt5a","Fzd7","sCs", 2,0.6666666666666666,0.310422,0.9312659999999999,0.02889080212084161,0.02889080212084161,3,1,4.478761333333333,13.436284,0.2317261188999595,0.2317261188999595,1.390306050616,12.512754455544,0.006694753447369345,0.006694753447369346)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $rowVals = $data[$i]
  for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item($r, $c).Value2 = $rowVals[$c - 1]
  }
}
